# =====================================================================
# Adds a new "2022-Q3" sheet (right after "总计") with fund-holding data,
# and shifts the summary ("总计") sheet rows down by one to make room for
# the new quarter, appending the old last row (2020-Q4) as a new row.
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right before "2022-Q2"
#    (duplicate "2021-Q4", which has the same layout/styles and enough
#    rows, so sheetPr/pageSetup/styling all come along for free).
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$template.Copy($beforeSheet)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q3"

# The template has 29 data rows; we only need 16 (header + 15 funds).
$newSheet.Range("A17:H29").EntireRow.Delete()

# Header row text
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows: A=index(n) B=基金代码 C=基金名称 D=基金规模 E=股票总仓位 F=仓位占比 G=持有市值(亿元) H=仓位排名
$data = @(
    @(0, '860001', '光大阳光混合A', '8.51', '90.45', '3.36', '0.2859', 10),
    @(1, '007040', '新疆前海联合泳隆灵活配置混合C', '6.78', '91.50', '3.49', '0.2366', 10),
    @(2, '860052', '光大阳光启明星创新驱动主题混合B', '1.63', '90.71', '4.23', '0.0689', 5),
    @(3, '860053', '光大阳光启明星创新驱动主题混合C', '1.21', '90.71', '4.23', '0.0512', 5),
    @(4, '004128', '新疆前海联合泳隆灵活配置混合A', '1.02', '91.50', '3.49', '0.0356', 10),
    @(5, '011351', '金鹰年年邮益一年持有期混合A', '3.43', '34.33', '0.81', '0.0278', 7),
    @(6, '519097', '新华中小市值优选混合', '0.71', '67.35', '3.14', '0.0223', 7),
    @(7, '860016', '光大阳光启明星创新驱动主题混合A', '0.42', '90.71', '4.23', '0.0178', 5),
    @(8, '000757', '华富智慧城市灵活配置混合', '0.50', '84.73', '2.82', '0.0141', 9),
    @(9, '005569', '中融智选红利股票A', '0.21', '92.04', '4.19', '0.0088', 4),
    @(10, '860036', '光大阳光混合B', '0.09', '90.45', '3.36', '0.0030', 10),
    @(11, '410006', '华富策略精选混合', '0.10', '70.87', '2.89', '0.0029', 5),
    @(12, '011352', '金鹰年年邮益一年持有期混合C', '0.27', '34.33', '0.81', '0.0022', 7),
    @(13, '005570', '中融智选红利股票C', '0.03', '92.04', '4.19', '0.0013', 4),
    @(14, '860037', '光大阳光混合C', '0.00', '90.45', '3.36', 0, 10)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Range("A" + $r).Value = $row[0]

    # Text columns that could otherwise be auto-converted to numbers:
    # force text storage, then drop the number-format style residue.
    $newSheet.Range("B" + $r).NumberFormat = "@"
    $newSheet.Range("B" + $r).Value = $row[1]
    $newSheet.Range("B" + $r).ClearFormats()

    $newSheet.Range("C" + $r).Value = $row[2]

    $newSheet.Range("D" + $r).NumberFormat = "@"
    $newSheet.Range("D" + $r).Value = $row[3]
    $newSheet.Range("D" + $r).ClearFormats()

    $newSheet.Range("E" + $r).NumberFormat = "@"
    $newSheet.Range("E" + $r).Value = $row[4]
    $newSheet.Range("E" + $r).ClearFormats()

    $newSheet.Range("F" + $r).NumberFormat = "@"
    $newSheet.Range("F" + $r).Value = $row[5]
    $newSheet.Range("F" + $r).ClearFormats()

    if ($r -eq 16) {
        # last row: G is a genuine numeric zero in the source data
        $newSheet.Range("G" + $r).Value = $row[6]
    } else {
        $newSheet.Range("G" + $r).NumberFormat = "@"
        $newSheet.Range("G" + $r).Value = $row[6]
        $newSheet.Range("G" + $r).ClearFormats()
    }

    $newSheet.Range("H" + $r).Value = $row[7]

    $r = $r + 1
}

$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: insert the 2022-Q3 row at the top
#    of the data, shifting every later row down by one, and materialize
#    a brand-new last row for the quarter that used to be last (2020-Q4).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("总计")

# Create row 9 using row 8 as a style template (keeps the s="2" styling
# on the index column), then fill in the values that used to sit in row 8.
$ws1.Range("A8:D8").Copy($ws1.Range("A9:D9"))
$ws1.Range("A9").Value = 7
$ws1.Range("B9").Value = $ws1.Range("B8").Value2
$ws1.Range("C9").Value = $ws1.Range("C8").Value2
$ws1.Range("D9").Value = $ws1.Range("D8").Value2

# Shift rows 7->8, 6->7, 5->6, 4->5, 3->4, 2->3 (bottom-up so sources
# are read before being overwritten).
$ws1.Range("B8").Value = $ws1.Range("B7").Value2
$ws1.Range("C8").Value = $ws1.Range("C7").Value2
$ws1.Range("D8").Value = $ws1.Range("D7").Value2

$ws1.Range("B7").Value = $ws1.Range("B6").Value2
$ws1.Range("C7").Value = $ws1.Range("C6").Value2
$ws1.Range("D7").Value = $ws1.Range("D6").Value2

$ws1.Range("B6").Value = $ws1.Range("B5").Value2
$ws1.Range("C6").Value = $ws1.Range("C5").Value2
$ws1.Range("D6").Value = $ws1.Range("D5").Value2

$ws1.Range("B5").Value = $ws1.Range("B4").Value2
$ws1.Range("C5").Value = $ws1.Range("C4").Value2
$ws1.Range("D5").Value = $ws1.Range("D4").Value2

$ws1.Range("B4").Value = $ws1.Range("B3").Value2
$ws1.Range("C4").Value = $ws1.Range("C3").Value2
$ws1.Range("D4").Value = $ws1.Range("D3").Value2

$ws1.Range("B3").Value = $ws1.Range("B2").Value2
$ws1.Range("C3").Value = $ws1.Range("C2").Value2
$ws1.Range("D3").Value = $ws1.Range("D2").Value2

# New 2022-Q3 figures go into row 2.
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 15
$ws1.Range("D2").Value = 0.78

$ws1.Range("A1").Select()

Write-Host "edit complete"
